# WS_holdings.xlsx edit
#
# 1) Cell A16 (shared-string footer) — bump the "as of" date from
#    2021-05-24 to 2021-05-25.
# 2) Weight / Percent Change table (D2:E13) — refresh the daily
#    figures to the new snapshot's values.
#
# The sheet ships with cell-level protection (password hash "D382"),
# so locked cells can't be written to until it's unprotected. We
# unprotect, make the edits, then re-protect the sheet afterward so it
# is left in a protected state again.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Unprotect()

# --- Footer disclosure text: update the "as of" date -----------------
$ws.Range("A16").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-25 for illustrative purposes only and are subject to change."

# --- Weight (D) / Percent Change (E) table, rows 2-13 -----------------
$newValues = @{
    2  = @{ D = 0.0310935635555849;  E = 0.0003374957813029322 }
    3  = @{ D = 0.02363555113263877; E = -0.01148886283704587 }
    4  = @{ D = 0.0521119951704264;  E = 0.001153934918070698 }
    5  = @{ D = 0.1380772260136681;  E = -0.002259522272433867 }
    6  = @{ D = 0.03148608431944987; E = -0.02100840336134435 }
    7  = @{ D = 0.1160224877020379;  E = 0.00206504904491478 }
    8  = @{ D = 0.1017368071047555;  E = -0.004028566196667138 }
    9  = @{ D = 0.02938139953873799; E = -0.01002865329512903 }
    10 = @{ D = 0.1267763067911162;  E = -0.01041666666666663 }
    11 = @{ D = 0.2462964622430162;  E = 0.0001819174094961351 }
    12 = @{ D = 0.1033821164285681;  E = -0.00134125311362332 }
    13 = @{ D = 0.9999999999999999; E = -0.003053738438972498 }
}

foreach ($row in $newValues.Keys) {
    $ws.Range("D$row").Value = $newValues[$row].D
    $ws.Range("E$row").Value = $newValues[$row].E
}

# --- Restore sheet protection -----------------------------------------
$ws.Protect()
